# Update cryptocurrency price/volume data (commit: "Updated cryptos list on Thu Jul 27 21:23:28 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.181.51"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "1.858.59"
$ws.Range("E3").Value = "  -1.21%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'0.7132"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").Value = "'240.27"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "'0.07731"
$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").Value = "'0.3075"
$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").Value = "'24.91"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "1.891.83"
$ws.Range("E12").Value = "  +1.37%  "

$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").Value = "'0.7141"
$ws.Range("E14").Value = "  -2.23%  "

$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "29.223.07"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "'5.865"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "'243.71"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "'0.000007805"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D21").Value = "2.115.19"
$ws.Range("E21").Value = "  -1.24%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").Value = "'7.931"
$ws.Range("E23").Value = "  +1.80%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").Value = "'0.1579"
$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").Value = "'162.63"
$ws.Range("E26").Value = "  -0.32%  "

$ws.Range("D27").Value = "'8.894"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("D28").Value = "'18.22"
$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").Value = "'1.321"
$ws.Range("E29").Value = "  -3.46%  "

$ws.Range("D30").Value = "'1.494"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("D31").Value = "'4.370"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").Value = "'4.121"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "'0.05184"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").Value = "'1.907"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").Value = "'0.7270"
$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("D37").Value = "'2.680"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("D39").Value = "'2.686"
$ws.Range("E39").Value = "  -1.15%  "

$ws.Range("D40").Value = "1.152.94"

$ws.Range("D41").Value = "'0.9000"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").Value = "'6.097"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").Value = "'72.18"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "'101.69"
$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("D46").Value = "2.014.77"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").Value = "'0.5234"
$ws.Range("E47").Value = "  -2.41%  "

$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").Value = "'9.283"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").Value = "'2.870"
$ws.Range("E51").Value = "  -0.37%  "
